# This script rotates the species-observation data (columns A,B,D,E,F,G,H,Q,R)
# down by one row across rows 5-9, with row 9's values wrapping around to row 5.
# (Column L, an empty placeholder cell, also follows the same rows as it is tied
# to the "Goodyera repens / VU" record.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for rows 5-9 or the columns that move together.
# NOTE: use .Value() (as a call) rather than the bare .Value property so the
# actual evaluated cell value is captured instead of a property reference.
$rows = 5..9
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        A = $ws.Cells.Item($r, 1).Value()
        B = $ws.Cells.Item($r, 2).Value()
        D = $ws.Cells.Item($r, 4).Value()
        E = $ws.Cells.Item($r, 5).Value()
        F = $ws.Cells.Item($r, 6).Value()
        G = $ws.Cells.Item($r, 7).Value()
        H = $ws.Cells.Item($r, 8).Value()
        L = $ws.Cells.Item($r, 12).Value()
        Q = $ws.Cells.Item($r, 17).Value()
        R = $ws.Cells.Item($r, 18).Value()
    }
}

# Write back the rotated values: row r gets the values that previously belonged
# to the row above it (row 5 wraps around and gets row 9's previous values).
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $prevIndex = ($i - 1 + $rows.Count) % $rows.Count
    $prevRow = $rows[$prevIndex]
    $src = $data[$prevRow]

    $ws.Cells.Item($r, 1).Value = $src.A
    $ws.Cells.Item($r, 2).Value = $src.B
    $ws.Cells.Item($r, 4).Value = $src.D
    $ws.Cells.Item($r, 5).Value = $src.E
    $ws.Cells.Item($r, 6).Value = $src.F
    $ws.Cells.Item($r, 7).Value = $src.G
    $ws.Cells.Item($r, 8).Value = $src.H
    if ($src.L -ne $null) {
        $ws.Cells.Item($r, 12).Value = $src.L
    } else {
        $ws.Cells.Item($r, 12).Value = ""
    }
    $ws.Cells.Item($r, 17).Value = $src.Q
    $ws.Cells.Item($r, 18).Value = $src.R
}
